# One-click update from Van Paper 07:23 AM on 2025-11-20
#
# Leaderboard update:
#  - Row 33 ("REBECCA'S BAKERY & CAFE") gets its Last Invoice Date filled in.
#  - Two brand-new prospects are inserted (keeping the sheet sorted by
#    Customer Number), pushing the two trailing rows down from 34/35 to 36/37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows at 34 and 35, shifting the old 34/35 down to 36/37 ---
# (done before touching row 33's formatting, so the new blank rows inherit
# row 33's original "blank date" look rather than the updated one)
$ws.Range("A34:A35").EntireRow.Insert()
$ws.Rows.Item(34).RowHeight = 13.05
$ws.Rows.Item(35).RowHeight = 13.05

# --- Update existing row 33: set its (previously blank) Last Invoice Date ---
# Pull the date format/alignment from an existing dated cell (xlPasteFormats)
# so we reuse the workbook's existing date style instead of growing new ones.
$ws.Range("D32").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value = 45980

# New row 34: SMALL HOURS LLC
$ws.Range("A34").Value = "SMALL HOURS LLC"
$ws.Range("B34").Value = "House Account"
$ws.Range("C34").Value = "030"
$ws.Range("E34").Value = "0008375"

# New row 35: YELLOW BRICK ROAD CHILD CARE
$ws.Range("A35").Value = "YELLOW BRICK ROAD CHILD CARE "
$ws.Range("B35").Value = "Dack, Suzanne"
$ws.Range("C35").Value = "023"
$ws.Range("E35").Value = "0008376"
